$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Summary
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 0.6554307116104869
$wsSummary.Range("C2").Value = 0.5992822966507177
$wsSummary.Range("D2").Value = 0.9382022471910112
$wsSummary.Range("E2").Value = 0.7313868613138687
$wsSummary.Range("F2").Value = 0.8428667563930013
$wsSummary.Range("G2").Value = 0.9182292400958691
$wsSummary.Range("H2").Value = 0.7780337780022164
$wsSummary.Range("I2").Value = 501
$wsSummary.Range("J2").Value = 335
$wsSummary.Range("K2").Value = 199
$wsSummary.Range("L2").Value = 33

# ---------------------------------------------------------------------------
# Sheet 2: Classification Report
# ---------------------------------------------------------------------------
$wsReport = $wb.Worksheets.Item("Classification Report")

$wsReport.Range("B2").Value = 0.8577586206896551
$wsReport.Range("C2").Value = 0.3726591760299626
$wsReport.Range("D2").Value = 0.5195822454308094

$wsReport.Range("B3").Value = 0.5992822966507177
$wsReport.Range("C3").Value = 0.9382022471910112
$wsReport.Range("D3").Value = 0.7313868613138687

$wsReport.Range("B4").Value = 0.6554307116104869
$wsReport.Range("C4").Value = 0.6554307116104869
$wsReport.Range("D4").Value = 0.6554307116104869
$wsReport.Range("E4").Value = 0.6554307116104869

$wsReport.Range("B5").Value = 0.7285204586701864
$wsReport.Range("C5").Value = 0.6554307116104869
$wsReport.Range("D5").Value = 0.625484553372339

$wsReport.Range("B6").Value = 0.7285204586701866
$wsReport.Range("C6").Value = 0.6554307116104869
$wsReport.Range("D6").Value = 0.625484553372339

# ---------------------------------------------------------------------------
# Sheet 3: Confusion Matrix
# ---------------------------------------------------------------------------
$wsConf = $wb.Worksheets.Item("Confusion Matrix")

$wsConf.Range("B2").Value = 199
$wsConf.Range("C2").Value = 335

$wsConf.Range("B3").Value = 33
$wsConf.Range("C3").Value = 501
